$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.799304008483887
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.812331199645996
$ws.Range("D1").Value = 1.307560682296753
$ws.Range("E1").Value = 0.9523842930793762
